$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.741.07"
$ws.Range("E2").Value = "  +2.17%  "

$ws.Range("D3").Value = "3.090.59"
$ws.Range("E3").Value = "  +5.37%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.06"
$ws.Range("E5").Value = "  +2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.20"
$ws.Range("E6").Value = "  +6.47%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.086.00"
$ws.Range("E8").Value = "  +5.38%  "

$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("E10").Value = "  -1.85%  "

$ws.Range("E11").Value = "  +3.65%  "

$ws.Range("E12").Value = "  +4.97%  "

$ws.Range("E13").Value = "  +2.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.39"
$ws.Range("E14").Value = "  +6.22%  "

$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").Value = "3.600.82"
$ws.Range("E16").Value = "  +5.37%  "

$ws.Range("D17").Value = "66.711.51"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("E18").Value = "  +3.25%  "

$ws.Range("D19").Value = "3.089.04"
$ws.Range("E19").Value = "  +5.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.29"
$ws.Range("E20").Value = "  +4.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.32"
$ws.Range("E21").Value = "  +5.18%  "

$ws.Range("E22").Value = "  +3.54%  "

$ws.Range("E23").Value = "  +3.47%  "

$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("E25").Value = "  +5.74%  "

$ws.Range("E26").Value = "  +7.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("E31").Value = "  +4.16%  "

$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.24"
$ws.Range("E33").Value = "  +4.47%  "

$ws.Range("E34").Value = "  +3.79%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  +3.73%  "

$ws.Range("E37").Value = "  +2.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.35"
$ws.Range("E38").Value = "  +5.88%  "

$ws.Range("E39").Value = "  +6.47%  "

$ws.Range("E40").Value = "  +6.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.30"
$ws.Range("E41").Value = "  +1.40%  "

$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.69"
$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("E45").Value = "  +3.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.57"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").Value = "2.786.10"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.08"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.93"
$ws.Range("E50").Value = "  +7.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  +1.83%  "
